# Locate the "twolineheader" worksheet (sheet5.xml) - it is the active sheet
# in the workbook (tabSelected=1), so ActiveSheet should already refer to it,
# but look it up explicitly by name to be safe.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("twolineheader")
$ws.Activate()

# Add a new header cell in B1 with text that contains special characters,
# and a numeric value below it in B2.
$ws.Range("B1").Value = "m'`"\ a"
$ws.Range("B2").Value = 3

# Update the selection to match the new last-edited cell.
$ws.Range("B2").Select()
